$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 6.240107999999999
$ws.Range("H2").Value = 18.720324
$ws.Range("I2").Value = 0.01732230523539376
$ws.Range("J2").Value = 0.01732230523539376
$ws.Range("M2").Value = 0.2901893333333334
$ws.Range("N2").Value = 0.870568
$ws.Range("O2").Value = 0.03429389578125064
$ws.Range("P2").Value = 0.03429389578125064
$ws.Range("Q2").Value = 1.810812780448
$ws.Range("R2").Value = 16.297315024032
$ws.Range("S2").Value = 0.0005940493304336059
$ws.Range("T2").Value = 0.0005940493304336061
$ws.Range("G3").Value = 6.240107999999999
$ws.Range("H3").Value = 18.720324
$ws.Range("I3").Value = 0.01732230523539376
$ws.Range("J3").Value = 0.01732230523539376
$ws.Range("O3").Value = 0.8402845891331153
$ws.Range("P3").Value = 0.8402845891331153
$ws.Range("Q3").Value = 44.36935607787199
$ws.Range("R3").Value = 399.3242047008479
$ws.Range("S3").Value = 0.01455566613756126
$ws.Range("T3").Value = 0.01455566613756126
$ws.Range("G4").Value = 6.240107999999999
$ws.Range("H4").Value = 18.720324
$ws.Range("I4").Value = 0.01732230523539376
$ws.Range("J4").Value = 0.01732230523539376
$ws.Range("O4").Value = 0.1254215150856341
$ws.Range("P4").Value = 0.1254215150856341
$ws.Range("Q4").Value = 6.622603740004
$ws.Range("R4").Value = 59.603433660036
$ws.Range("S4").Value = 0.002172589767398897
$ws.Range("T4").Value = 0.002172589767398897
$ws.Range("I5").Value = 0.9592798330716089
$ws.Range("J5").Value = 0.9592798330716091
$ws.Range("M5").Value = 0.2901893333333334
$ws.Range("N5").Value = 0.870568
$ws.Range("O5").Value = 0.03429389578125064
$ws.Range("P5").Value = 0.03429389578125064
$ws.Range("Q5").Value = 100.2797351822907
$ws.Range("R5").Value = 902.517616640616
$ws.Range("S5").Value = 0.03289744262041327
$ws.Range("T5").Value = 0.03289744262041328
$ws.Range("I6").Value = 0.9592798330716089
$ws.Range("J6").Value = 0.9592798330716091
$ws.Range("O6").Value = 0.8402845891331153
$ws.Range("P6").Value = 0.8402845891331153
$ws.Range("S6").Value = 0.8060680603962603
$ws.Range("T6").Value = 0.8060680603962606
$ws.Range("I7").Value = 0.9592798330716089
$ws.Range("J7").Value = 0.9592798330716091
$ws.Range("O7").Value = 0.1254215150856341
$ws.Range("P7").Value = 0.1254215150856341
$ws.Range("Q7").Value = 366.7485432152437
$ws.Range("R7").Value = 3300.736888937193
$ws.Range("S7").Value = 0.1203143300549354
$ws.Range("T7").Value = 0.1203143300549354
$ws.Range("G8").Value = 8.428738666666666
$ws.Range("I8").Value = 0.02339786169299727
$ws.Range("J8").Value = 0.02339786169299728
$ws.Range("M8").Value = 0.2901893333333334
$ws.Range("N8").Value = 0.870568
$ws.Range("O8").Value = 0.03429389578125064
$ws.Range("P8").Value = 0.03429389578125064
$ws.Range("Q8").Value = 2.445930054520889
$ws.Range("R8").Value = 22.013370490688
$ws.Range("S8").Value = 0.000802403830403765
$ws.Range("T8").Value = 0.0008024038304037652
$ws.Range("G9").Value = 8.428738666666666
$ws.Range("I9").Value = 0.02339786169299727
$ws.Range("J9").Value = 0.02339786169299728
$ws.Range("O9").Value = 0.8402845891331153
$ws.Range("P9").Value = 0.8402845891331153
$ws.Range("Q9").Value = 59.93128759769244
$ws.Range("R9").Value = 539.381588379232
$ws.Range("S9").Value = 0.01966086259929367
$ws.Range("T9").Value = 0.01966086259929367
$ws.Range("G10").Value = 8.428738666666666
$ws.Range("I10").Value = 0.02339786169299727
$ws.Range("J10").Value = 0.02339786169299728
$ws.Range("O10").Value = 0.1254215150856341
$ws.Range("P10").Value = 0.1254215150856341
$ws.Range("S10").Value = 0.002934595263299838
$ws.Range("T10").Value = 0.002934595263299838

Write-Output "Updated 83 cells"
